# Bursdager.xlsx — remove the two test rows ("Donald Duck - Test" / "Mikke
# Mus - Test") that were appended at the bottom of the birthday list, and
# leave the now-empty date cells (still styled) behind, matching the
# upload diff. Also reposition the view/selection the way Excel left it
# after the edit (scrolled down, A54:C55 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 54 was: A54="Donald", B54="Duck - Test", C54=45980 (date, style 4)
# Row 55 was: A55="Mikke",  B55="Mus - Test",  C55=45981 (date, style 4)
# Target: A54/B54/A55/B55 removed entirely; C54/C55 cleared but keep style.
$ws.Range("A54:B55").ClearContents()
$ws.Range("C54:C55").ClearContents()

# Scroll the view down and select A54:C55, as in the saved workbook.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$ws.Range("A54:C55").Select()
